$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35
$ws.Range("A35:J35").NumberFormat = "@"
$ws.Range("L35:N35").NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "990210"
$ws.Cells.Item($row, 2).Value = "Facilitator"
$ws.Cells.Item($row, 3).Value = "aa"
$ws.Cells.Item($row, 4).Value = "aa"
$ws.Cells.Item($row, 5).Value = "a@gmail.com"
$ws.Cells.Item($row, 6).Value = "78908908908"
$ws.Cells.Item($row, 7).Value = "gujkhjkyhjghk"
$ws.Cells.Item($row, 8).Value = "Muslim"
$ws.Cells.Item($row, 9).Value = "Female"
$ws.Cells.Item($row, 10).Value = "Married"
$ws.Cells.Item($row, 11).Value = 77
$ws.Cells.Item($row, 12).Value = "no"
$ws.Cells.Item($row, 13).Value = "hhhhhhhhhhhh"
$ws.Cells.Item($row, 14).Value = "hhhhhhhh"

# Reset cell style back to the default (no explicit style) to match the
# unstyled cells used by the other data rows in the sheet.
$ws.Range("A35:N35").Style = "Normal"
